$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 13960.23
$ws.Range("I125").Value = 16085.375
$ws.Range("J125").Value = 10560
$ws.Range("K125").Value = 144768.375
$ws.Range("L125").Value = 95040
$ws.Range("M125").Value = -142308.375
$ws.Range("N125").Value = -99960

$ws.Range("H132").Value = 5297.8477
$ws.Range("I132").Value = 1019.4286
$ws.Range("J132").Value = 50221.25
$ws.Range("K132").Value = 3058.2858
$ws.Range("L132").Value = 150663.75
$ws.Range("M132").Value = -528.2857999999997
$ws.Range("N132").Value = -155723.75

$ws.Range("H135").Value = 765.76666
$ws.Range("I135").Value = 652.7406999999999
$ws.Range("J135").Value = 1783
$ws.Range("K135").Value = 5874.6663
$ws.Range("L135").Value = 16047
$ws.Range("M135").Value = -3339.6663
$ws.Range("N135").Value = -21117

$ws.Range("H138").Value = 4287.7393
$ws.Range("I138").Value = 1044.7778
$ws.Range("J138").Value = 6372.5
$ws.Range("K138").Value = 3134.3334
$ws.Range("L138").Value = 19117.5
$ws.Range("M138").Value = 2005.6666
$ws.Range("N138").Value = -29397.5

$ws.Range("H141").Value = 2541.9788
$ws.Range("I141").Value = 2277.8667
$ws.Range("J141").Value = 8484.5
$ws.Range("K141").Value = 6833.6001
$ws.Range("L141").Value = 25453.5
$ws.Range("M141").Value = -1653.6001
$ws.Range("N141").Value = -35813.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4908.909
$ws.Range("I2").Value = 3571.4285
$ws.Range("J2").Value = 7249.5
$ws.Range("K2").Value = 3571.4285
$ws.Range("L2").Value = 7249.5
$ws.Range("M2").Value = -3458.4285
$ws.Range("N2").Value = -7475.5

$ws.Range("H22").Value = 8969.6
$ws.Range("I22").Value = 4972
$ws.Range("J22").Value = 9969
$ws.Range("K22").Value = 4972
$ws.Range("L22").Value = 9969
$ws.Range("M22").Value = -4673
$ws.Range("N22").Value = -10567

$ws.Range("H32").Value = 8456.940000000001
$ws.Range("I32").Value = 7684.6
$ws.Range("J32").Value = 15408
$ws.Range("K32").Value = 7684.6
$ws.Range("L32").Value = 15408
$ws.Range("M32").Value = -7397.6
$ws.Range("N32").Value = -15982

$ws.Range("H61").Value = 19238366
$ws.Range("I61").Value = 19238366
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 19238366
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -19238154

$ws.Range("H63").Value = 5296.143
$ws.Range("I63").Value = 2874.75
$ws.Range("J63").Value = 8524.666999999999
$ws.Range("K63").Value = 2874.75
$ws.Range("L63").Value = 8524.666999999999
$ws.Range("M63").Value = -2188.75
$ws.Range("N63").Value = -9896.666999999999

$ws.Range("H66").Value = 5296.143
$ws.Range("I66").Value = 2874.75
$ws.Range("J66").Value = 8524.666999999999
$ws.Range("K66").Value = 14373.75
$ws.Range("L66").Value = 42623.335
$ws.Range("M66").Value = -10941.75
$ws.Range("N66").Value = -49487.335

$ws.Range("H74").Value = 3116.5
$ws.Range("I74").Value = 1987.8334
$ws.Range("J74").Value = 5656
$ws.Range("K74").Value = 1987.8334
$ws.Range("L74").Value = 5656
$ws.Range("M74").Value = -1113.8334
$ws.Range("N74").Value = -7404

$ws.Range("H77").Value = 3116.5
$ws.Range("I77").Value = 1987.8334
$ws.Range("J77").Value = 5656
$ws.Range("K77").Value = 9939.166999999999
$ws.Range("L77").Value = 28280
$ws.Range("M77").Value = -5571.166999999999
$ws.Range("N77").Value = -37016

$ws.Range("H116").Value = 4908.909
$ws.Range("I116").Value = 3571.4285
$ws.Range("J116").Value = 7249.5
$ws.Range("K116").Value = 3571.4285
$ws.Range("L116").Value = 7249.5
$ws.Range("M116").Value = -1277.4285
$ws.Range("N116").Value = -11837.5

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 4592.4634
$ws.Range("I132").Value = 4130.5
$ws.Range("J132").Value = 7918.6
$ws.Range("K132").Value = 12391.5
$ws.Range("L132").Value = 23755.8
$ws.Range("M132").Value = -9861.5
$ws.Range("N132").Value = -28815.8

$ws.Range("H136").Value = 19238366
$ws.Range("I136").Value = 19238366
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 57715098
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -57712548

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4908.909
$ws.Range("I3").Value = 3571.4285
$ws.Range("J3").Value = 7249.5
$ws.Range("K3").Value = 3571.4285
$ws.Range("L3").Value = 7249.5
$ws.Range("M3").Value = -3457.4285
$ws.Range("N3").Value = -7477.5

$ws.Range("H11").Value = 188.3077
$ws.Range("I11").Value = 36.25
$ws.Range("J11").Value = 431.6
$ws.Range("K11").Value = 36.25
$ws.Range("L11").Value = 431.6
$ws.Range("M11").Value = 103.75
$ws.Range("N11").Value = -711.6

$ws.Range("H86").Value = 952.625
$ws.Range("I86").Value = 934.8182
$ws.Range("J86").Value = 1036.5714
$ws.Range("K86").Value = 934.8182
$ws.Range("L86").Value = 1036.5714
$ws.Range("M86").Value = 188.1818
$ws.Range("N86").Value = -3282.5714

$ws.Range("H89").Value = 952.625
$ws.Range("I89").Value = 934.8182
$ws.Range("J89").Value = 1036.5714
$ws.Range("K89").Value = 4674.091
$ws.Range("L89").Value = 5182.857
$ws.Range("M89").Value = 941.9089999999997
$ws.Range("N89").Value = -16414.857

$ws.Range("H94").Value = 4340
$ws.Range("I94").Value = 600
$ws.Range("J94").Value = 5275
$ws.Range("K94").Value = 600
$ws.Range("L94").Value = 5275
$ws.Range("M94").Value = -149
$ws.Range("N94").Value = -6177

$ws.Range("H134").Value = 4768.5127
$ws.Range("I134").Value = 5057.5557
$ws.Range("J134").Value = 1300
$ws.Range("K134").Value = 15172.6671
$ws.Range("L134").Value = 3900
$ws.Range("M134").Value = -12637.6671
$ws.Range("N134").Value = -8970

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5561.212
$ws.Range("I31").Value = 2671.0667
$ws.Range("J31").Value = 7969.6665
$ws.Range("K31").Value = 2671.0667
$ws.Range("L31").Value = 7969.6665
$ws.Range("M31").Value = -2376.0667
$ws.Range("N31").Value = -8559.666499999999

$ws.Range("H34").Value = 5561.212
$ws.Range("I34").Value = 2671.0667
$ws.Range("J34").Value = 7969.6665
$ws.Range("K34").Value = 2671.0667
$ws.Range("L34").Value = 7969.6665
$ws.Range("M34").Value = -2469.0667
$ws.Range("N34").Value = -8373.666499999999

$ws.Range("H58").Value = 6684.4165
$ws.Range("I58").Value = 5521.75
$ws.Range("J58").Value = 7265.75
$ws.Range("K58").Value = 5521.75
$ws.Range("L58").Value = 7265.75
$ws.Range("M58").Value = -5318.75
$ws.Range("N58").Value = -7671.75

$ws.Range("H62").Value = 8799.200000000001
$ws.Range("I62").Value = 8499
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 8499
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -7875
$ws.Range("N62").Value = -11248

$ws.Range("H65").Value = 8799.200000000001
$ws.Range("I65").Value = 8499
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 42495
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -39375
$ws.Range("N65").Value = -56240

$ws.Range("H99").Value = 7717.864
$ws.Range("I99").Value = 6906.857
$ws.Range("J99").Value = 9137.125
$ws.Range("K99").Value = 6906.857
$ws.Range("L99").Value = 9137.125
$ws.Range("M99").Value = -5408.857
$ws.Range("N99").Value = -12133.125

$ws.Range("H107").Value = 1262.5883
$ws.Range("I107").Value = 1402.1818
$ws.Range("J107").Value = 1006.6667
$ws.Range("K107").Value = 1402.1818
$ws.Range("L107").Value = 1006.6667
$ws.Range("M107").Value = 517.8181999999999
$ws.Range("N107").Value = -4846.6667

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H126").Value = 7717.864
$ws.Range("I126").Value = 6906.857
$ws.Range("J126").Value = 9137.125
$ws.Range("K126").Value = 20720.571
$ws.Range("L126").Value = 27411.375
$ws.Range("M126").Value = -18250.571
$ws.Range("N126").Value = -32351.375

$ws.Range("H136").Value = 6684.4165
$ws.Range("I136").Value = 5521.75
$ws.Range("J136").Value = 7265.75
$ws.Range("K136").Value = 16565.25
$ws.Range("L136").Value = 21797.25
$ws.Range("M136").Value = -14015.25
$ws.Range("N136").Value = -26897.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 574.25
$ws.Range("I103").Value = 408.8
$ws.Range("J103").Value = 850
$ws.Range("K103").Value = 1226.4
$ws.Range("L103").Value = 2550
$ws.Range("M103").Value = -347.4000000000001
$ws.Range("N103").Value = -4308

$ws.Range("H117").Value = 910301.8
$ws.Range("I117").Value = 2046
$ws.Range("J117").Value = 1250897.8
$ws.Range("K117").Value = 6138
$ws.Range("L117").Value = 3752693.4
$ws.Range("M117").Value = -2696
$ws.Range("N117").Value = -3759577.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11102.608
$ws.Range("I70").Value = 9595.200000000001
$ws.Range("J70").Value = 12262.154
$ws.Range("K70").Value = 9595.200000000001
$ws.Range("L70").Value = 12262.154
$ws.Range("M70").Value = -9325.200000000001
$ws.Range("N70").Value = -12802.154

$ws.Range("H73").Value = 11102.608
$ws.Range("I73").Value = 9595.200000000001
$ws.Range("J73").Value = 12262.154
$ws.Range("K73").Value = 9595.200000000001
$ws.Range("L73").Value = 12262.154
$ws.Range("M73").Value = -8659.200000000001
$ws.Range("N73").Value = -14134.154

$ws.Range("H80").Value = 3866.4119
$ws.Range("I80").Value = 3282.75
$ws.Range("J80").Value = 4385.222
$ws.Range("K80").Value = 3282.75
$ws.Range("L80").Value = 4385.222
$ws.Range("M80").Value = -2284.75
$ws.Range("N80").Value = -6381.222

$ws.Range("H83").Value = 3866.4119
$ws.Range("I83").Value = 3282.75
$ws.Range("J83").Value = 4385.222
$ws.Range("K83").Value = 16413.75
$ws.Range("L83").Value = 21926.11
$ws.Range("M83").Value = -11421.75
$ws.Range("N83").Value = -31910.11

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 538.2727
$ws.Range("I55").Value = 690.75
$ws.Range("J55").Value = 131.66667
$ws.Range("K55").Value = 690.75
$ws.Range("L55").Value = 131.66667
$ws.Range("M55").Value = -517.75
$ws.Range("N55").Value = -477.66667

$ws.Range("H68").Value = 4991.25
$ws.Range("I68").Value = 2819.6
$ws.Range("J68").Value = 8610.666999999999
$ws.Range("K68").Value = 2819.6
$ws.Range("L68").Value = 8610.666999999999
$ws.Range("M68").Value = -2070.6
$ws.Range("N68").Value = -10108.667

$ws.Range("H71").Value = 4991.25
$ws.Range("I71").Value = 2819.6
$ws.Range("J71").Value = 8610.666999999999
$ws.Range("K71").Value = 14098
$ws.Range("L71").Value = 43053.335
$ws.Range("M71").Value = -10354
$ws.Range("N71").Value = -50541.335

$ws.Range("H82").Value = 2227.1
$ws.Range("I82").Value = 1880.8334
$ws.Range("J82").Value = 2746.5
$ws.Range("K82").Value = 1880.8334
$ws.Range("L82").Value = 2746.5
$ws.Range("M82").Value = -1519.8334
$ws.Range("N82").Value = -3468.5

$ws.Range("H85").Value = 2227.1
$ws.Range("I85").Value = 1880.8334
$ws.Range("J85").Value = 2746.5
$ws.Range("K85").Value = 1880.8334
$ws.Range("L85").Value = 2746.5
$ws.Range("M85").Value = -632.8334
$ws.Range("N85").Value = -5242.5

$ws.Range("H93").Value = 8778.909
$ws.Range("I93").Value = 2722.9092
$ws.Range("J93").Value = 20890.908
$ws.Range("K93").Value = 2722.9092
$ws.Range("L93").Value = 20890.908
$ws.Range("M93").Value = -1474.9092
$ws.Range("N93").Value = -23386.908

$ws.Range("H132").Value = 8420.68
$ws.Range("I132").Value = 9209
$ws.Range("J132").Value = 3986.375
$ws.Range("K132").Value = 27627
$ws.Range("L132").Value = 11959.125
$ws.Range("M132").Value = -25097
$ws.Range("N132").Value = -17019.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 7000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 7000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 7000
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -7280

$ws.Range("H136").Value = 1741.6613
$ws.Range("I136").Value = 1080.4036
$ws.Range("J136").Value = 9280
$ws.Range("K136").Value = 3241.2108
$ws.Range("L136").Value = 27840
$ws.Range("M136").Value = -691.2108000000003
$ws.Range("N136").Value = -32940
